# BOT; UPDATE DATA
# Appends the latest day's row (2020-05-16, serial 43967) to the three
# data sheets ("all", "kobe", "other"), pushing the trailing footnote
# row down by one, and tweaks a few already-reported figures on "kobe"
# for the previous day (2020-05-15, serial 43966) now that more data
# came in.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "all": new row 39 (data), footnote moves from row 39 -> 40
# ---------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Activate()

# Push the trailing "footnote" row down one row first, carrying its
# formatting along so the label keeps its original look.
$wsAll.Range("B39").Copy()
$wsAll.Range("B40").PasteSpecial(-4122)

# New data row: copy the format of the previous data row (38) down
# into row 39, then overwrite with the new day's figures.
$wsAll.Range("A38:H38").Copy()
$wsAll.Range("A39:H39").PasteSpecial(-4122)
$wsAll.Range("A39").Value = 43967
$wsAll.Range("B39").Value = 282
$wsAll.Range("C39").Value = 278
$wsAll.Range("D39").Value = 63
$wsAll.Range("E39").Value = 53
$wsAll.Range("F39").Value = 10
$wsAll.Range("G39").Value = 11
$wsAll.Range("H39").Value = 204

# Restore the footnote text on its new row (reuses the existing
# shared string already used elsewhere on this sheet).
$wsAll.Range("B40").Value = "※　24・34・53・58・59・60・158・161・163・192・237・248・268・272・276・277例目（計16件）は市外在住者です。レイメケイケンシガイザイジュウシャ"

$wsAll.Range("A39").Select()

# ---------------------------------------------------------------
# Sheet "kobe": revise row 93 (2020-05-15) and append new row 94
# (2020-05-16); footnote moves from row 94 -> 95
# ---------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Activate()

# Revised counts for the already-recorded day.
$wsKobe.Range("D93").Value = 1
$wsKobe.Range("E93").Value = 282
$wsKobe.Range("F93").Value = 58

# Push the trailing "footnote" row down one row first.
$wsKobe.Range("B94").Copy()
$wsKobe.Range("B95").PasteSpecial(-4122)

# New data row, formatted like the row above it.
$wsKobe.Range("A93:J93").Copy()
$wsKobe.Range("A94:J94").PasteSpecial(-4122)
$wsKobe.Range("A94").Value = 43967
$wsKobe.Range("B94").Value = 0
$wsKobe.Range("C94").Value = 2813
$wsKobe.Range("D94").Value = 0
$wsKobe.Range("E94").Value = 282
$wsKobe.Range("F94").Value = 58
$wsKobe.Range("G94").Value = 49
$wsKobe.Range("H94").Value = 9
$wsKobe.Range("I94").Value = 11
$wsKobe.Range("J94").Value = 195

$wsKobe.Range("B95").Value = "※　24・34・53・58・59・60・158・161・163・192・237・248・268・272・276・277例目（計16件）は市外在住者です。レイメケイケンシガイザイジュウシャ"

$wsKobe.Range("A94").Select()

# ---------------------------------------------------------------
# Sheet "other": new row 69 (data), footnote moves from row 69 -> 70
# ---------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Activate()

$wsOther.Range("B69").Copy()
$wsOther.Range("B70").PasteSpecial(-4122)

$wsOther.Range("A68:H68").Copy()
$wsOther.Range("A69:H69").PasteSpecial(-4122)
$wsOther.Range("A69").Value = 43967
$wsOther.Range("B69").Value = 0
$wsOther.Range("C69").Value = 14
$wsOther.Range("D69").Value = 5
$wsOther.Range("E69").Value = 4
$wsOther.Range("F69").Value = 1
$wsOther.Range("G69").Value = 0
$wsOther.Range("H69").Value = 9

$wsOther.Range("B70").Value = "※他自治体において、3月10日以前の感染者の発生はございません。タジチタイニチ"

$wsOther.Range("A68").Select()

# ---------------------------------------------------------------
# Leave "all" as the active sheet/tab, matching the source workbook.
# ---------------------------------------------------------------
$wsAll.Activate()
